$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 196, shifting existing rows 196:234 down to 197:235
$ws.Rows.Item(196).Insert()

# Fill in the new row 196 with the new data point
$ws.Range("A196").Value = 3
$ws.Range("B196").Value = "Femacal de La Calera"
$ws.Range("C196").Value = "Coquimbo"
$ws.Range("D196").Value = 44504
$ws.Range("E196").Value = 5
$ws.Range("F196").Value = 100112040
$ws.Range("G196").Value = "Cilantro"
$ws.Range("H196").Value = "Sin especificar"
$ws.Range("I196").Value = "Primera"
$ws.Range("J196").Value = 168
$ws.Range("K196").Value = 2500
$ws.Range("L196").Value = 2500
$ws.Range("M196").Value = 2500
$ws.Range("N196").Value = "`$/docena de atados (3 kilos)"
$ws.Range("O196").Value = "Provincia de Quillota"
$ws.Range("P196").Value = 833
$ws.Range("Q196").Value = 3
$ws.Range("R196").Value = "Hortaliza"
